# Applies the "Tirusse" worksheet update:
#  - adds a new "Подбор" label in column D for rows 6 and 16
#  - lowers the packet-length inputs on rows 15/16 (32->7, 45->25)
#  - turns the B21 "Real packet duration" cell from a formula (=B19) into a
#    fixed/overridden value of 6.8, which cascades through the dependent
#    formulas in B23 and B24
#  - moves the view/selection down to B24

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 6: tag the second measurement block's "Service info" row with D6 = "Подбор"
$ws.Range("D6").Value = "Подбор"

# Row 15: Packet length 32 -> 7 bytes
$ws.Range("B15").Value = 7

# Row 16: Service info 45 -> 25 bytes, and tag it with D16 = "Подбор" as well
$ws.Range("B16").Value = 25
$ws.Range("D16").Value = "Подбор"

# Row 21: replace the formula "=B19" with a fixed, experimentally measured
# value of 6.8 (downstream formulas in B23/B24 recalculate automatically)
$ws.Range("B21").Value = 6.8

# Scroll/selection: bring row 6 to the top of the view and select B24
$ws.Range("B24").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
